$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at the top of the data (rows 6-8), pushing existing
# rows 6-18 down to rows 9-21 (and consequently 16-18 to 19-21).
$ws.Rows("6:8").Insert()

# Fill the three newly inserted rows with the new weekly price records.

# Row 6: Especial
$ws.Cells.Item(6,1).Value = 3
$ws.Cells.Item(6,2).Value = "Femacal de La Calera"
$ws.Cells.Item(6,3).Value = "Coquimbo"
$ws.Cells.Item(6,4).Value = 44452
$ws.Cells.Item(6,5).Value = 5
$ws.Cells.Item(6,6).Value = "Fruta"
$ws.Cells.Item(6,7).Value = 100107
$ws.Cells.Item(6,8).Value = "Otros"
$ws.Cells.Item(6,9).Value = 100107002
$ws.Cells.Item(6,10).Value = "Chirimoya"
$ws.Cells.Item(6,11).Value = "Cultivar IV Región"
$ws.Cells.Item(6,12).Value = "Especial"
$ws.Cells.Item(6,13).Value = 56
$ws.Cells.Item(6,14).Value = 30000
$ws.Cells.Item(6,15).Value = 30000
$ws.Cells.Item(6,16).Value = 30000
$ws.Cells.Item(6,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(6,18).Value = "Provincia del Elquí"
$ws.Cells.Item(6,19).Value = 3000
$ws.Cells.Item(6,20).Value = 10

# Row 7: Primera
$ws.Cells.Item(7,1).Value = 3
$ws.Cells.Item(7,2).Value = "Femacal de La Calera"
$ws.Cells.Item(7,3).Value = "Coquimbo"
$ws.Cells.Item(7,4).Value = 44452
$ws.Cells.Item(7,5).Value = 5
$ws.Cells.Item(7,6).Value = "Fruta"
$ws.Cells.Item(7,7).Value = 100107
$ws.Cells.Item(7,8).Value = "Otros"
$ws.Cells.Item(7,9).Value = 100107002
$ws.Cells.Item(7,10).Value = "Chirimoya"
$ws.Cells.Item(7,11).Value = "Cultivar IV Región"
$ws.Cells.Item(7,12).Value = "Primera"
$ws.Cells.Item(7,13).Value = 60
$ws.Cells.Item(7,14).Value = 27000
$ws.Cells.Item(7,15).Value = 27000
$ws.Cells.Item(7,16).Value = 27000
$ws.Cells.Item(7,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(7,18).Value = "Provincia del Elquí"
$ws.Cells.Item(7,19).Value = 2700
$ws.Cells.Item(7,20).Value = 10

# Row 8: Segunda
$ws.Cells.Item(8,1).Value = 3
$ws.Cells.Item(8,2).Value = "Femacal de La Calera"
$ws.Cells.Item(8,3).Value = "Coquimbo"
$ws.Cells.Item(8,4).Value = 44452
$ws.Cells.Item(8,5).Value = 5
$ws.Cells.Item(8,6).Value = "Fruta"
$ws.Cells.Item(8,7).Value = 100107
$ws.Cells.Item(8,8).Value = "Otros"
$ws.Cells.Item(8,9).Value = 100107002
$ws.Cells.Item(8,10).Value = "Chirimoya"
$ws.Cells.Item(8,11).Value = "Cultivar IV Región"
$ws.Cells.Item(8,12).Value = "Segunda"
$ws.Cells.Item(8,13).Value = 50
$ws.Cells.Item(8,14).Value = 25000
$ws.Cells.Item(8,15).Value = 25000
$ws.Cells.Item(8,16).Value = 25000
$ws.Cells.Item(8,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8,18).Value = "Provincia del Elquí"
$ws.Cells.Item(8,19).Value = 2500
$ws.Cells.Item(8,20).Value = 10
